# Daily attendance processing - 2025-10-02 16:11:56
# Applies updated attendance / coverage statistics and "Recorded By" ordering
# changes to the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain "General number format" reference cell (style index matching the
# existing percentage-text cells). Used below to make sure the percentage
# strings we write keep being stored as plain text (not auto-converted to a
# numeric percentage) while keeping the original cell style untouched.
$plainFormat = $ws.Range("L8")

function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $plainFormat.Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats
}

# ---------------------------------------------------------------------------
# Class Statistics block (K/L columns)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 40        # Recorded Sessions
$ws.Range("L7").Value = 9         # Missing Sessions
Set-TextValue $ws.Range("L9") "26.1%"    # Coverage %
Set-TextValue $ws.Range("L10") "48.5%"   # Average Attendance %

# ---------------------------------------------------------------------------
# Group Statistics block (O/P/R/S columns)
# ---------------------------------------------------------------------------
# Row 19 - Year 2 / B1
$ws.Range("O19").Value = 5
$ws.Range("P19").Value = 0
Set-TextValue $ws.Range("R19") "29.4%"
Set-TextValue $ws.Range("S19") "59.7%"

# Row 22 - Year 2 / B4
$ws.Range("O22").Value = 4
$ws.Range("P22").Value = 1
Set-TextValue $ws.Range("R22") "23.5%"
Set-TextValue $ws.Range("S22") "24.9%"

# ---------------------------------------------------------------------------
# "Recorded By" (column G) reordering on already-recorded sessions
# ---------------------------------------------------------------------------
$ws.Range("G14").Value = "marian.samir@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"
$ws.Range("G17").Value = "ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G31").Value = "marian.samir@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"
$ws.Range("G34").Value = "ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G45").Value = "System, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, backup@backdoor.com"
$ws.Range("G51").Value = "Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G62").Value = "System, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, backup@backdoor.com"
$ws.Range("G68").Value = "Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G72").Value = "mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G77").Value = "mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G83").Value = "Youstina.ibrahim@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G85").Value = "wafaa.ebida@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("G98").Value = "user@user.com, nourhanmohamed@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G102").Value = "wafaa.ebida@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("G116").Value = "enas.omran@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G119").Value = "Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, marinasorial@med.asu.edu.eg"
$ws.Range("G133").Value = "enas.omran@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G136").Value = "Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, marinasorial@med.asu.edu.eg"
$ws.Range("G149").Value = "Youstina.ibrahim@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# Sessions that just got recorded: row 81 (Year 2 / B1) and row 131 (Year 2 / B4)
# Copy the formatting from an already-"Recorded" row (row 7) so the
# pink "Not Recorded" fill becomes the green "Recorded" fill, then fill
# in the reviewer, student count and status.
# ---------------------------------------------------------------------------
$recordedFormat = $ws.Range("A7:I7")

$recordedFormat.Copy()
$ws.Range("A81:I81").PasteSpecial(-4122)
$ws.Range("G81").Value = "mariam.youssif.std@med.asu.edu.eg"
$ws.Range("H81").Value = "112/154"
$ws.Range("I81").Value = "Recorded"

$recordedFormat.Copy()
$ws.Range("A131:I131").PasteSpecial(-4122)
$ws.Range("G131").Value = "mariam.youssif.std@med.asu.edu.eg"
$ws.Range("H131").Value = "50/226"
$ws.Range("I131").Value = "Recorded"

$excel.CutCopyMode = 0
